$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date-serial value that must be bumped
# from 46061 to 46062 (one day later) for every data row (2 through 486).
for ($r = 2; $r -le 486; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value = 46062
    }
}
